$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new log entry row (row 26), copying the date cell's format from
# the previous row so it keeps the date-style formatting.
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A26").Value = 45656

$ws.Range("B26").Value = "Fleshing out area more with intereactable objects etc."

$ws.Range("C26").Value = 6

# Update the selection to reflect the new active cell
$ws.Range("A26").Select()
